$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(3,  3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(4,  1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642),
    @(5,  0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795),
    @(6,  0.1554434735375247, 0.05231270169004087, 3.082599426703578, 0.4998867070740569, 3.790242309005201),
    @(7,  3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 12.0302756157461),
    @(8,  0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.896700893398075),
    @(9,  0.3464964993005633, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 10.24810883928323),
    @(10, 3.182878228561681, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 4.733082622659194),
    @(11, 1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245),
    @(12, 3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729),
    @(13, 0.000009318123435519965, 0.004309184025731883, 0.1529057820181812, 0.4998867070740569, 0.6571109912414055),
    @(14, 3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126),
    @(15, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(16, 3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(17, 1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
